# Insert two new weekly price rows for "Coliflor" (Terminal La Palmera de
# La Serena) right before the existing row 863, shifting the remainder of
# the dataset (old rows 863-961) down by two rows to 865-963.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 863 (pushes rows 863.. down to 865..)
$ws.Rows.Item(863).Insert()
$ws.Rows.Item(863).Insert()

# Constant columns shared by every data row in this sheet.
$mercadoId = 8
$mercado   = "Terminal La Palmera de La Serena"
$region    = "Coquimbo"
$codreg    = 4
$catId     = 100112008
$categoria = "Coliflor"
$variedad  = "Sin especificar"
$unidad    = "`$/unidad"
$origen    = "Provincia del Elqu$([char]0xED)"
$kgUnid    = 1
$clasif    = "Hortaliza"

# --- New row 863 ("Primera") ---
$ws.Range("A863").Value = $mercadoId
$ws.Range("B863").Value = $mercado
$ws.Range("C863").Value = $region
$ws.Range("D863").Value = 44946
$ws.Range("E863").Value = $codreg
$ws.Range("F863").Value = $catId
$ws.Range("G863").Value = $categoria
$ws.Range("H863").Value = $variedad
$ws.Range("I863").Value = "Primera"
$ws.Range("J863").Value = 2400
$ws.Range("K863").Value = 800
$ws.Range("L863").Value = 900
$ws.Range("M863").Value = 850
$ws.Range("N863").Value = $unidad
$ws.Range("O863").Value = $origen
$ws.Range("P863").Value = 850
$ws.Range("Q863").Value = $kgUnid
$ws.Range("R863").Value = $clasif

# --- New row 864 ("Segunda") ---
$ws.Range("A864").Value = $mercadoId
$ws.Range("B864").Value = $mercado
$ws.Range("C864").Value = $region
$ws.Range("D864").Value = 44946
$ws.Range("E864").Value = $codreg
$ws.Range("F864").Value = $catId
$ws.Range("G864").Value = $categoria
$ws.Range("H864").Value = $variedad
$ws.Range("I864").Value = "Segunda"
$ws.Range("J864").Value = 1460
$ws.Range("K864").Value = 600
$ws.Range("L864").Value = 700
$ws.Range("M864").Value = 650
$ws.Range("N864").Value = $unidad
$ws.Range("O864").Value = $origen
$ws.Range("P864").Value = 650
$ws.Range("Q864").Value = $kgUnid
$ws.Range("R864").Value = $clasif

Write-Output ("New dimension: " + $ws.UsedRange.Address())
